$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "AZIONE" mapping column (H) is being folded into the "STATI UNITI D'EUROPA"
# mapping column (G), which is then relabeled "CENTRO". Any row that had a value
# in both G and H gets the sum in G; any row with a value only in H gets that
# value moved into G. Afterwards column H is removed entirely, shifting I:P left
# to H:O.

for ($r = 2; $r -le 58; $r++) {
    $gCell = $ws.Cells.Item($r, 7)   # column G
    $hCell = $ws.Cells.Item($r, 8)   # column H

    $gText = $gCell.Text
    $hText = $hCell.Text

    $gHasValue = -not [string]::IsNullOrEmpty($gText)
    $hHasValue = -not [string]::IsNullOrEmpty($hText)

    if ($hHasValue) {
        if ($gHasValue) {
            $gCell.Value = [double]$gText + [double]$hText
        } else {
            $gCell.Value = [double]$hText
        }
        $hCell.Value = $null
    }
}

# Remove the now-redundant column H; I:P shift left to H:O.
$ws.Columns("H").Delete()

# Relabel the merged column header as "CENTRO".
$ws.Cells.Item(1, 7).Value = "CENTRO"

# Restore the active selection recorded in the saved workbook.
$ws.Range("F3").Select()
